# Adds two new columns, I ("I0") and J ("IF"), to the data table on the
# active sheet, mirroring the header style used by the existing columns
# and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the formatting of the last existing header cell (H1, bold/centered
# with borders) onto the two new header cells before writing their text,
# so I1/J1 pick up the same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows ----------------------------------------------------------
# For every data row (2-38), column I is 1 and column J repeats the value
# already present in column H - except row 3, which carries different
# values (I3=4, J3=6).
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()

    if ($r -eq 3) {
        $ws.Cells.Item($r, 9).Value = 4
        $ws.Cells.Item($r, 10).Value = 6
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
